$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 11445.185
$ws.Range("J87").Value = 11445.185
$ws.Range("L87").Value = 11445.185
$ws.Range("N87").Value = -13941.185

$ws.Range("H90").Value = 11445.185
$ws.Range("J90").Value = 11445.185
$ws.Range("L90").Value = 34335.555
$ws.Range("N90").Value = -46815.555

$ws.Range("H107").Value = 893
$ws.Range("I107").Value = 916.1177
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 916.1177
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1003.8823
$ws.Range("N107").Value = -4340

$ws.Range("H129").Value = 938.3333
$ws.Range("I129").Value = 367.5
$ws.Range("J129").Value = 1395
$ws.Range("K129").Value = 1102.5
$ws.Range("L129").Value = 4185
$ws.Range("M129").Value = 3897.5
$ws.Range("N129").Value = -14185

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3204.7727
$ws.Range("I32").Value = 2619.2856
$ws.Range("J32").Value = 15500
$ws.Range("K32").Value = 2619.2856
$ws.Range("L32").Value = 15500
$ws.Range("M32").Value = -2332.2856
$ws.Range("N32").Value = -16074

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("N62").Value = 0
$ws.Range("L62").Value = ""

$ws.Range("H63").Value = 3698.2942
$ws.Range("I63").Value = 2919.3572
$ws.Range("J63").Value = 7333.3335
$ws.Range("K63").Value = 2919.3572
$ws.Range("L63").Value = 7333.3335
$ws.Range("M63").Value = -2233.3572
$ws.Range("N63").Value = -8705.333500000001

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("N64").Value = 0
$ws.Range("L64").Value = ""

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("N65").Value = 0
$ws.Range("L65").Value = ""

$ws.Range("H66").Value = 3698.2942
$ws.Range("I66").Value = 2919.3572
$ws.Range("J66").Value = 7333.3335
$ws.Range("K66").Value = 14596.786
$ws.Range("L66").Value = 36666.6675
$ws.Range("M66").Value = -11164.786
$ws.Range("N66").Value = -43530.6675

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("N67").Value = 0
$ws.Range("L67").Value = ""

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("N68").Value = 0
$ws.Range("L68").Value = ""
$ws.Range("M68").Value = ""

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("N71").Value = 0
$ws.Range("L71").Value = ""
$ws.Range("M71").Value = ""

$ws.Range("H76").Value = 8000
$ws.Range("J76").Value = 8000
$ws.Range("L76").Value = 8000
$ws.Range("N76").Value = -8676

$ws.Range("H79").Value = 8000
$ws.Range("J79").Value = 8000
$ws.Range("L79").Value = 8000
$ws.Range("N79").Value = -10340

$ws.Range("H132").Value = 1072.3572
$ws.Range("I132").Value = 582.42255
$ws.Range("J132").Value = 3748.1538
$ws.Range("K132").Value = 1747.26765
$ws.Range("L132").Value = 11244.4614
$ws.Range("M132").Value = 782.73235
$ws.Range("N132").Value = -16304.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("N62").Value = 0
$ws.Range("L62").Value = ""

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("N65").Value = 0
$ws.Range("L65").Value = ""

$ws.Range("H75").Value = 10868.5
$ws.Range("I75").Value = 10868.5
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 10868.5
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = -9932.5
$ws.Range("M75").Value = ""

$ws.Range("H76").Value = 18314
$ws.Range("J76").Value = 18314
$ws.Range("L76").Value = 18314
$ws.Range("N76").Value = -18944

$ws.Range("H78").Value = 10868.5
$ws.Range("I78").Value = 10868.5
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 32605.5
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = -27925.5
$ws.Range("M78").Value = ""

$ws.Range("H79").Value = 18314
$ws.Range("J79").Value = 18314
$ws.Range("L79").Value = 18314
$ws.Range("N79").Value = -20498

$ws.Range("H82").Value = 16881.334
$ws.Range("I82").Value = 6333.3335
$ws.Range("K82").Value = 6333.3335
$ws.Range("M82").Value = -5950.3335

$ws.Range("H85").Value = 16881.334
$ws.Range("I85").Value = 6333.3335
$ws.Range("K85").Value = 6333.3335
$ws.Range("M85").Value = -5007.3335

$ws.Range("H86").Value = 28575640
$ws.Range("I86").Value = 47621456
$ws.Range("J86").Value = 6915.7856
$ws.Range("K86").Value = 47621456
$ws.Range("L86").Value = 6915.7856
$ws.Range("M86").Value = -47620333
$ws.Range("N86").Value = -9161.785599999999

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("N88").Value = 0
$ws.Range("L88").Value = ""

$ws.Range("H89").Value = 28575640
$ws.Range("I89").Value = 47621456
$ws.Range("J89").Value = 6915.7856
$ws.Range("K89").Value = 238107280
$ws.Range("L89").Value = 34578.928
$ws.Range("M89").Value = -238101664
$ws.Range("N89").Value = -45810.928

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("N91").Value = 0
$ws.Range("L91").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 18010.5
$ws.Range("J74").Value = 18010.5
$ws.Range("L74").Value = 18010.5
$ws.Range("N74").Value = -19758.5

$ws.Range("H77").Value = 18010.5
$ws.Range("J77").Value = 18010.5
$ws.Range("L77").Value = 54031.5
$ws.Range("N77").Value = -62767.5

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("N82").Value = 0
$ws.Range("L82").Value = ""

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("N85").Value = 0
$ws.Range("L85").Value = ""

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("N87").Value = 0
$ws.Range("L87").Value = ""
$ws.Range("M87").Value = ""

$ws.Range("H88").Value = 20000
$ws.Range("J88").Value = 20000
$ws.Range("L88").Value = 20000
$ws.Range("N88").Value = -20812

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("N90").Value = 0
$ws.Range("L90").Value = ""
$ws.Range("M90").Value = ""

$ws.Range("H91").Value = 20000
$ws.Range("J91").Value = 20000
$ws.Range("L91").Value = 20000
$ws.Range("N91").Value = -22808

$ws.Range("H132").Value = 33839.613
$ws.Range("I132").Value = 765.0909
$ws.Range("J132").Value = 114688.445
$ws.Range("K132").Value = 2295.2727
$ws.Range("L132").Value = 344065.335
$ws.Range("M132").Value = 234.7273
$ws.Range("N132").Value = -349125.335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15626161
$ws.Range("I131").Value = 460
$ws.Range("J131").Value = 16667874
$ws.Range("K131").Value = 1380
$ws.Range("L131").Value = 50003622
$ws.Range("M131").Value = 3660
$ws.Range("N131").Value = -50013702

$ws.Range("H139").Value = 5171
$ws.Range("I139").Value = 6065.5557
$ws.Range("J139").Value = 3158.25
$ws.Range("K139").Value = 18196.6671
$ws.Range("L139").Value = 9474.75
$ws.Range("M139").Value = -13056.6671
$ws.Range("N139").Value = -19754.75

$ws.Range("H141").Value = 8342.727999999999
$ws.Range("I141").Value = 9596.25
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 28788.75
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -23608.75
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 8422
$ws.Range("J95").Value = 8422
$ws.Range("L95").Value = 8422
$ws.Range("N95").Value = -13914

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12050.5625
$ws.Range("I62").Value = 6434
$ws.Range("J62").Value = 19271.857
$ws.Range("K62").Value = 6434
$ws.Range("L62").Value = 19271.857
$ws.Range("M62").Value = -5810
$ws.Range("N62").Value = -20519.857

$ws.Range("H65").Value = 12050.5625
$ws.Range("I65").Value = 6434
$ws.Range("J65").Value = 19271.857
$ws.Range("K65").Value = 32170
$ws.Range("L65").Value = 96359.285
$ws.Range("M65").Value = -29050
$ws.Range("N65").Value = -102599.285

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("N97").Value = 0
$ws.Range("L97").Value = ""

$ws.Range("H132").Value = 20356260
$ws.Range("I132").Value = 29763134
$ws.Range("J132").Value = 1542509.8
$ws.Range("K132").Value = 89289402
$ws.Range("L132").Value = 4627529.4
$ws.Range("M132").Value = -89286872
$ws.Range("N132").Value = -4632589.4

$ws.Range("H139").Value = 48800
$ws.Range("J139").Value = 48800
$ws.Range("L139").Value = 48800
$ws.Range("N139").Value = -59080

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("N141").Value = 0
$ws.Range("L141").Value = ""
